$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 17, pushing rows 17-25 down to 19-27
$ws.Range("A17:C18").EntireRow.Insert()

# Fill in the two new rows (16 and 17) with the new time entries
$ws.Range("A16").Value = "Statistical Analysis (centering, CI, ICC)"
$ws.Range("B16").Value = "May. 20"
$ws.Range("C16").Value = 3

$ws.Range("A17").Value = "Statistical Analysis (centering, CI, ICC)"
$ws.Range("B17").Value = "May. 21"
$ws.Range("C17").Value = 6

# Update the TOTAL formula (now at row 25) to include the new rows, sum C2:C23
$ws.Range("C25").Formula = "=SUM(C2:C23)"

# Update the Amount label (now at row 27) and formula
$ws.Range("A27").Value = "Amount (25`$/hour)"
$ws.Range("C27").Formula = "=C25*25"

# Update selection to match target (E28)
$ws.Range("E28").Select()
